$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.682599065821893, 1.6902474096035132, -1.3475849219033682, 2.6638145849986126)"
$ws.Range("C2").Value = "NIG(1.973811874207114, 1.6449359498414733, 2.469656053739959, 5.622544467228233)"
$ws.Range("D2").Value = "JSU(-0.741201265403811, 1.0060574798467825, 0.679260594338078, 2.7259575353294823)"
$ws.Range("E2").Value = "JSU(-0.7714322775378564, 1.0536051183745592, 4.734279630206152, 4.103008216991709)"
